$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Professional summary paragraph: "all Black and Asian-American voters"
#    -> "50M voters" (plain text swap, single run paragraph).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
  "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Discovered systematic demographic coding errors affecting 50M voters, developed", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Siege Analytics bullet: split "... affecting all Black and
#    Asian-American voters, developed ..." so "50M" becomes its own bold run,
#    matching the formatting already used for the other bolded stats.
#    Scope the Find to the specific paragraph so the unique-match quirk of
#    Content.Find (which can mutate the wrong occurrence when the same text
#    appears more than once in the story) can't bite us.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $para = $d.Paragraphs($i)
  if ($para.Range.Text -like "*Discovered systematic race coding errors affecting all Black and Asian-American voters*") {
    $bulletRange = $para.Range
    $bulletRange.Find.Execute("all Black and Asian-American") | Out-Null
    $bulletRange.Text = "50M"
    $bulletRange.Bold = 1
    $bulletRange.Font.Color = 5258796
    break
  }
}

# ---------------------------------------------------------------------------
# 3) Reorder PROFESSIONAL EXPERIENCE entries:
#      Siege Analytics, Mautinoa, Salsa Labs, Praxis Project, PCCC, Helm
#    becomes
#      Siege Analytics, Helm, Mautinoa, PCCC, Salsa Labs, Praxis Project
#
#    Every one of these five job blocks is exactly 5 paragraphs long with the
#    identical paragraph-style pattern (Heading3, Normal, Normal, Normal,
#    Normal), so instead of physically inserting/deleting paragraphs (which
#    risks corrupting styles/rsids) we simply capture the text/formatting of
#    each of those 25 paragraphs and write it back into the same 25
#    paragraph "slots" in the new order.
# ---------------------------------------------------------------------------

function Get-ParaIndexByText($marker) {
  for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like $marker) {
      return $i
    }
  }
  return -1
}

$mautinoaStart = Get-ParaIndexByText "Software Engineer - Mautinoa Technologies*"
$salsaStart    = Get-ParaIndexByText "Software Engineer - Salsa Labs*"
$praxisStart   = Get-ParaIndexByText "Interim Technology Manager - The Praxis Project*"
$pcccStart     = Get-ParaIndexByText "Research Director - PCCC*"
$helmStart     = Get-ParaIndexByText "Data Products Manager - Helm*"

# Each block is 5 paragraphs: job title (Heading3), subtitle, and 3 bullets.
$mautinoaTexts = @()
$salsaTexts    = @()
$praxisTexts   = @()
$pcccTexts     = @()
$helmTexts     = @()
# Paragraph.Range.Text includes the trailing paragraph-mark (CR, chr 13);
# strip it so writing the captured string back doesn't insert a spurious
# extra paragraph break.
for ($k = 0; $k -lt 5; $k++) {
  $mautinoaTexts += , $d.Paragraphs($mautinoaStart + $k).Range.Text.TrimEnd([char]13)
  $salsaTexts    += , $d.Paragraphs($salsaStart + $k).Range.Text.TrimEnd([char]13)
  $praxisTexts   += , $d.Paragraphs($praxisStart + $k).Range.Text.TrimEnd([char]13)
  $pcccTexts     += , $d.Paragraphs($pcccStart + $k).Range.Text.TrimEnd([char]13)
  $helmTexts     += , $d.Paragraphs($helmStart + $k).Range.Text.TrimEnd([char]13)
}

# New order of blocks starting at $mautinoaStart (the first of the five slots):
$newBlocks = @($helmTexts, $mautinoaTexts, $pcccTexts, $salsaTexts, $praxisTexts)

$slot = $mautinoaStart
foreach ($block in $newBlocks) {
  foreach ($line in $block) {
    $d.Paragraphs($slot).Range.Text = $line
    $slot = $slot + 1
  }
}

# The "Modernized legacy ETL ... 57%" bullet (last line of the Helm block)
# carries a bold "57%" run in the source; Range.Text above flattened it back
# to plain text, so re-apply the bold/colored run now that it lives in its
# new slot.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $para = $d.Paragraphs($i)
  if ($para.Range.Text -like "*Modernized legacy ETL processes*57%*") {
    $bulletRange = $para.Range
    $bulletRange.Find.Execute("57%") | Out-Null
    $bulletRange.Text = "57%"
    $bulletRange.Bold = 1
    $bulletRange.Font.Color = 5258796
    break
  }
}
